$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Bold the section / sub-section heading paragraphs.
#    (bolding the whole paragraph Range also marks the paragraph-mark's
#    run properties as bold, matching <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>)
# ---------------------------------------------------------------------
$headings = @(
    "I. METİN NEŞRİNDE NÜSHALAR (tespit, seçim, tavsif)",
    "I. 1. Nüshaların tespit ve temini",
    "I. 2. Nüshaların seçimi",
    "I. 3. Nüsha Tavsifi",
    "II. TENKİTLİ METİN KURULMASINDA TEMEL HUSUSLAR",
    "II. 1. Tenkitli metin niçin kurulur?",
    "II. 2. Varyant tercihinde kriterler",
    "II. 3. Naşirin metne müdahalesi meselesi",
    "III. OKUYUCU KILAVUZU: YÖNTEM AÇIKLAMASI",
    "IV. DİĞER HUSUSLAR"
)

foreach ($h in $headings) {
    $rng = $d.Content
    $found = $rng.Find.Execute($h, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Paragraphs(1).Range.Font.Bold = 1
    }
}

# ---------------------------------------------------------------------
# 2. Join "istisnasız bütün nüshalar temin" + bookmark + " edilmelidir."
#    into a single sentence "istisnasız bütün nüshalar temin edilmelidir."
#    (removes the old _GoBack bookmark from this spot as a side effect).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("istisnasız bütün nüshalar temin edilmelidir.", $true, $false, $false, $false, $false, $true, 1, $false, "istisnasız bütün nüshalar temin edilmelidir.", 2)

# ---------------------------------------------------------------------
# 3. Relocate the (hidden, Word-managed) "_GoBack" bookmark to sit in the
#    middle of "Nüsha ailelerini" -> "Nüsha ailelerin" | "i ...", matching
#    the last-edit location left behind after editing the "Nüsha
#    ailelerini oluşturmada ..." sentence. Adding a bookmark with a name
#    that already exists elsewhere simply relocates it (names are unique).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nüsha ailelerin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
